$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 3.312999963760376
$ws.Range("E2").Value = 2.877000093460083
$ws.Range("F2").Value = 3.799000024795532
$ws.Range("G2").Value = 2.834000110626221
$ws.Range("H2").Value = 594273308
$ws.Range("I2").Value = "SMCI"

$ws.Range("D3").Value = 2.980000019073486
$ws.Range("E3").Value = 2.66700005531311
$ws.Range("F3").Value = 3.046999931335449
$ws.Range("G3").Value = 2.459000110626221
$ws.Range("H3").Value = 594273308
$ws.Range("I3").Value = "SMCI"

$ws.Range("D4").Value = 2.720999956130981
$ws.Range("E4").Value = 2.821000099182129
$ws.Range("F4").Value = 3.209000110626221
$ws.Range("G4").Value = 2.510999917984009
$ws.Range("H4").Value = 594273308
$ws.Range("I4").Value = "SMCI"

$ws.Range("D5").Value = 2.388000011444092
$ws.Range("E5").Value = 2.977999925613404
$ws.Range("F5").Value = 3.160000085830688
$ws.Range("G5").Value = 2.125
$ws.Range("H5").Value = 594273308
$ws.Range("I5").Value = "SMCI"

$ws.Range("D6").Value = 3.36299991607666
$ws.Range("E6").Value = 2.690999984741211
$ws.Range("F6").Value = 3.469000101089477
$ws.Range("G6").Value = 2.609999895095825
$ws.Range("H6").Value = 594273308
$ws.Range("I6").Value = "SMCI"

$ws.Range("D7").Value = 2.476999998092652
$ws.Range("E7").Value = 2.154999971389771
$ws.Range("F7").Value = 2.655999898910522
$ws.Range("G7").Value = 1.860000014305115
$ws.Range("H7").Value = 594273308
$ws.Range("I7").Value = "SMCI"

$ws.Range("D8").Value = 2.336999893188477
$ws.Range("E8").Value = 2.369999885559082
$ws.Range("F8").Value = 2.490000009536743
$ws.Range("G8").Value = 2.125999927520752
$ws.Range("H8").Value = 594273308
$ws.Range("I8").Value = "SMCI"

$ws.Range("D9").Value = 2.825000047683716
$ws.Range("E9").Value = 2.644999980926514
$ws.Range("F9").Value = 3.174999952316284
$ws.Range("G9").Value = 2.484999895095825
$ws.Range("H9").Value = 594273308
$ws.Range("I9").Value = "SMCI"

$ws.Range("D10").Value = 2.539999961853028
$ws.Range("E10").Value = 2.440000057220459
$ws.Range("F10").Value = 2.575000047683716
$ws.Range("G10").Value = 2.230000019073486
$ws.Range("H10").Value = 594273308
$ws.Range("I10").Value = "SMCI"

$ws.Range("D11").Value = 2.480000019073486
$ws.Range("E11").Value = 2.684999942779541
$ws.Range("F11").Value = 2.759999990463257
$ws.Range("G11").Value = 2.384999990463257
$ws.Range("H11").Value = 594273308
$ws.Range("I11").Value = "SMCI"

$ws.Range("D12").Value = 2.234999895095825
$ws.Range("E12").Value = 1.990000009536743
$ws.Range("F12").Value = 2.349999904632568
$ws.Range("G12").Value = 1.769999980926514
$ws.Range("H12").Value = 594273308
$ws.Range("I12").Value = "SMCI"

$ws.Range("D13").Value = 2.089999914169312
$ws.Range("E13").Value = 2.282999992370605
$ws.Range("F13").Value = 2.480000019073486
$ws.Range("G13").Value = 2.085000038146973
$ws.Range("H13").Value = 594273308
$ws.Range("I13").Value = "SMCI"

$ws.Range("D14").Value = 1.700000047683716
$ws.Range("E14").Value = 1.769999980926514
$ws.Range("F14").Value = 1.840000033378601
$ws.Range("G14").Value = 1.625
$ws.Range("H14").Value = 594273308
$ws.Range("I14").Value = "SMCI"

$ws.Range("D15").Value = 2.359999895095825
$ws.Range("E15").Value = 2.210000038146973
$ws.Range("F15").Value = 2.494999885559082
$ws.Range("G15").Value = 2.167999982833862
$ws.Range("H15").Value = 594273308
$ws.Range("I15").Value = "SMCI"

$ws.Range("D16").Value = 2.069999933242798
$ws.Range("E16").Value = 1.309999942779541
$ws.Range("F16").Value = 2.160000085830688
$ws.Range("G16").Value = 0.8500000238418579
$ws.Range("H16").Value = 594273308
$ws.Range("I16").Value = "SMCI"

$ws.Range("D17").Value = 1.378999948501587
$ws.Range("E17").Value = 1.509999990463257
$ws.Range("F17").Value = 1.585000038146973
$ws.Range("G17").Value = 1.378999948501587
$ws.Range("H17").Value = 594273308
$ws.Range("I17").Value = "SMCI"

$ws.Range("D18").Value = 2.105000019073486
$ws.Range("E18").Value = 2.242000102996826
$ws.Range("F18").Value = 2.243000030517578
$ws.Range("G18").Value = 2.105000019073486
$ws.Range("H18").Value = 594273308
$ws.Range("I18").Value = "SMCI"

$ws.Range("D19").Value = 1.932999968528748
$ws.Range("E19").Value = 1.830999970436096
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 1.830000042915344
$ws.Range("H19").Value = 594273308
$ws.Range("I19").Value = "SMCI"

$ws.Range("D20").Value = 1.919999957084656
$ws.Range("E20").Value = 2.068000078201294
$ws.Range("F20").Value = 2.144000053405762
$ws.Range("G20").Value = 1.855000019073486
$ws.Range("H20").Value = 594273308
$ws.Range("I20").Value = "SMCI"

$ws.Range("D21").Value = 2.384999990463257
$ws.Range("E21").Value = 2.796000003814697
$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 2.299999952316284
$ws.Range("H21").Value = 594273308
$ws.Range("I21").Value = "SMCI"

$ws.Range("D22").Value = 2.046000003814697
$ws.Range("E22").Value = 2.289999961853028
$ws.Range("F22").Value = 2.392999887466431
$ws.Range("G22").Value = 1.947999954223633
$ws.Range("H22").Value = 594273308
$ws.Range("I22").Value = "SMCI"

$ws.Range("D23").Value = 2.838000059127808
$ws.Range("E23").Value = 3.030999898910522
$ws.Range("F23").Value = 3.030999898910522
$ws.Range("G23").Value = 2.658999919891357
$ws.Range("H23").Value = 594273308
$ws.Range("I23").Value = "SMCI"

$ws.Range("D24").Value = 2.658999919891357
$ws.Range("E24").Value = 2.272000074386597
$ws.Range("F24").Value = 2.819000005722046
$ws.Range("G24").Value = 2.194999933242798
$ws.Range("H24").Value = 594273308
$ws.Range("I24").Value = "SMCI"

$ws.Range("D25").Value = 3.180000066757202
$ws.Range("E25").Value = 3.099999904632568
$ws.Range("F25").Value = 3.507999897003174
$ws.Range("G25").Value = 3.052000045776367
$ws.Range("H25").Value = 594273308
$ws.Range("I25").Value = "SMCI"

$ws.Range("D26").Value = 3.933000087738037
$ws.Range("E26").Value = 3.70199990272522
$ws.Range("F26").Value = 4.071000099182129
$ws.Range("G26").Value = 3.690999984741211
$ws.Range("H26").Value = 594273308
$ws.Range("I26").Value = "SMCI"

$ws.Range("D27").Value = 3.535000085830688
$ws.Range("E27").Value = 3.803999900817871
$ws.Range("F27").Value = 3.815999984741211
$ws.Range("G27").Value = 3.288000106811523
$ws.Range("H27").Value = 594273308
$ws.Range("I27").Value = "SMCI"

$ws.Range("D28").Value = 3.671999931335449
$ws.Range("E28").Value = 3.539000034332275
$ws.Range("F28").Value = 3.779999971389771
$ws.Range("G28").Value = 3.477999925613404
$ws.Range("H28").Value = 594273308
$ws.Range("I28").Value = "SMCI"

$ws.Range("D29").Value = 4.460999965667725
$ws.Range("E29").Value = 4.052000045776367
$ws.Range("F29").Value = 4.798999786376953
$ws.Range("G29").Value = 3.815999984741211
$ws.Range("H29").Value = 594273308
$ws.Range("I29").Value = "SMCI"

$ws.Range("D30").Value = 3.829999923706055
$ws.Range("E30").Value = 4.210000038146973
$ws.Range("F30").Value = 4.684000015258789
$ws.Range("G30").Value = 3.41100001335144
$ws.Range("H30").Value = 594273308
$ws.Range("I30").Value = "SMCI"

$ws.Range("D31").Value = 4.007999897003174
$ws.Range("E31").Value = 5.401000022888184
$ws.Range("F31").Value = 5.514999866485596
$ws.Range("G31").Value = 3.700999975204468
$ws.Range("H31").Value = 594273308
$ws.Range("I31").Value = "SMCI"

$ws.Range("D32").Value = 5.65500020980835
$ws.Range("E32").Value = 6.959000110626221
$ws.Range("F32").Value = 7.019000053405762
$ws.Range("G32").Value = 5.027999877929688
$ws.Range("H32").Value = 594273308
$ws.Range("I32").Value = "SMCI"

$ws.Range("D33").Value = 8.307999610900879
$ws.Range("E33").Value = 7.232999801635742
$ws.Range("F33").Value = 9.279999732971191
$ws.Range("G33").Value = 6.901999950408936
$ws.Range("H33").Value = 594273308
$ws.Range("I33").Value = "SMCI"

$ws.Range("D34").Value = 10.6569995880127
$ws.Range("E34").Value = 10.54300022125244
$ws.Range("F34").Value = 11.92399978637695
$ws.Range("G34").Value = 9.319000244140623
$ws.Range("H34").Value = 594273308
$ws.Range("I34").Value = "SMCI"

$ws.Range("D35").Value = 25.48900032043457
$ws.Range("E35").Value = 33.02700042724609
$ws.Range("F35").Value = 33.77899932861328
$ws.Range("G35").Value = 24.51000022888184
$ws.Range("H35").Value = 594273308
$ws.Range("I35").Value = "SMCI"

$ws.Range("D36").Value = 27.5
$ws.Range("E36").Value = 23.94700050354004
$ws.Range("F36").Value = 31.75
$ws.Range("G36").Value = 22.65900039672852
$ws.Range("H36").Value = 594273308
$ws.Range("I36").Value = "SMCI"

$ws.Range("D37").Value = 28
$ws.Range("E37").Value = 52.96099853515625
$ws.Range("F37").Value = 55.44400024414063
$ws.Range("G37").Value = 27.58799934387207
$ws.Range("H37").Value = 594273308
$ws.Range("I37").Value = "SMCI"

$ws.Range("D38").Value = 101
$ws.Range("E38").Value = 85.87999725341797
$ws.Range("F38").Value = 106.9000015258789
$ws.Range("G38").Value = 67.09999847412109
$ws.Range("H38").Value = 594273308
$ws.Range("I38").Value = "SMCI"

$ws.Range("D39").Value = 83.14399719238281
$ws.Range("E39").Value = 70.16500091552734
$ws.Range("F39").Value = 96.33000183105467
$ws.Range("G39").Value = 65.63099670410156
$ws.Range("H39").Value = 594273308
$ws.Range("I39").Value = "SMCI"

$ws.Range("D40").Value = 41.75
$ws.Range("E40").Value = 29.11000061035156
$ws.Range("F40").Value = 50.61000061035156
$ws.Range("G40").Value = 27.21999931335449
$ws.Range("H40").Value = 594273308
$ws.Range("I40").Value = "SMCI"

$ws.Range("D41").Value = 30.96999931335449
$ws.Range("E41").Value = 28.52000045776367
$ws.Range("F41").Value = 38.5
$ws.Range("G41").Value = 27.21999931335449
$ws.Range("H41").Value = 594273308
$ws.Range("I41").Value = "SMCI"

$ws.Range("D42").Value = 34.25
$ws.Range("E42").Value = 31.86000061035156
$ws.Range("F42").Value = 37.81999969482422
$ws.Range("G42").Value = 27.60000038146973
$ws.Range("H42").Value = 594273308
$ws.Range("I42").Value = "SMCI"

$ws.Range("D43").Value = 48.2400016784668
$ws.Range("E43").Value = 58.97000122070312
$ws.Range("F43").Value = 62.36000061035156
$ws.Range("G43").Value = 46.22000122070312
$ws.Range("H43").Value = 594273308
$ws.Range("I43").Value = "SMCI"
